$d = $word.ActiveDocument
$lb = [char]11

# The document content got "rotated" through a sequence of paragraphs:
# the paragraph formatting / headings / bold labels stay exactly where
# they are, but the body text that follows them shifts to a different
# paragraph. We replace, paragraph-by-paragraph (scoped Range so we
# never touch the wrong occurrence), the old body text with the new
# body text.

# Para 6: "Objetivos" (PT) body
$p = $d.Paragraphs.Item(6)
$p.Range.Text = "Sistemas de vácuo. Criogenia e baixa temperatura."

# Para 7: "Objetivos" (EN, italic) body
$p = $d.Paragraphs.Item(7)
$p.Range.Text = "Vacuum systems. Cryogenics and low temperature."

# Para 9: "Docente(s) Responsável(eis)" list item
$p = $d.Paragraphs.Item(9)
$p.Range.Text = "Fornecer os conhecimentos sobre sistemas de vácuo e técnicas de produção e utilização de baixas temperaturas."

# Para 11: "Programa resumido" (PT) body
$p = $d.Paragraphs.Item(11)
$p.Range.Text = "Teoria dos gases rarefeitos. Escoamento de gases. Bombas de vácuo. Descrição quantitativa do bombeamento de sistemas de vácuo. Medidores de pressão. Acessórios: armadilhas, anteparos, válvulas, etc. Adsorção, dessorção e evaporação de moléculas em vácuo. Detecção de vazamento. Vedação. Soldagem. Limpeza. " + $lb + "Criogenia. Propriedades de gases e líquidos criogênicos. Métodos para obtenção de baixa temperatura. Liquefação de gases. Medição de temperatura. Componentes criogênicos. Cálculo de transferência de calor em criostatos e dewars."

# Para 12: "Programa resumido" (EN, italic) body
$p = $d.Paragraphs.Item(12)
$p.Range.Text = "Provide knowledge about vacuum systems and production techniques and use of low temperatures."

# Para 14: "Programa" body (two lines joined by a w:br collapse into one
# line; include the line-break char in the search so the w:br is
# consumed too and no stray xml:space="preserve" is introduced)
$p = $d.Paragraphs.Item(14)
$r = $p.Range
$oldText = "Teoria dos gases rarefeitos. Escoamento de gases. Bombas de vácuo. Descrição quantitativa do bombeamento de sistemas de vácuo. Medidores de pressão. Acessórios: armadilhas, anteparos, válvulas, etc. Adsorção, dessorção e evaporação de moléculas em vácuo. Detecção de vazamento. Vedação. Soldagem. Limpeza. " + $lb + "Criogenia. Propriedades de gases e líquidos criogênicos. Métodos para obtenção de baixa temperatura. Liquefação de gases. Medição de temperatura. Componentes criogênicos. Cálculo de transferência de calor em criostatos e dewars."
$r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "Experimentos desenvolvidos em laboratório didático, realização de relatórios para cada experimento e de testes sobre o experimento em estudo.", 2) | Out-Null

# Para 17: "Avaliação" list item (Método:/Critério:/Norma de recuperação:)
# Only the (non-bold) value runs after each bold label change; the bold
# label runs are left completely untouched. Replace in reverse order
# (Norma -> Critério -> Método) so that no intermediate search text is
# ever duplicated within the paragraph before it is located/replaced.
$p = $d.Paragraphs.Item(17)

$biblio = "ROTH, A. Vacuum Technology, North-Holland, 1990." + $lb + "HARRIS, N. S. Modern Vacuum Practice, McGraw-Hill, 1989." + $lb + "HABLANIAN, M. H. High-Vacuum Technology, Marcel Dekker, 1997." + $lb + "BARRON, R. F. Cryogenic Systems, Oxford University Press, 1985." + $lb + "WEISEND, J. G. The Handbook of Cryogenic Engineering, Boca Raton: CRC Press, 1998."
$r = $p.Range
$r.Find.Execute("Aplicação de uma prova escrita e prática dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação", $true, $false, $false, $false, $false, $true, 1, $false, $biblio, 2) | Out-Null

$r = $p.Range
$r.Find.Execute("Média aritmética de duas provas escritas, testes, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + P2 + TR)/3", $true, $false, $false, $false, $false, $true, 1, $false, "Aplicação de uma prova escrita e prática dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação", 2) | Out-Null

$r = $p.Range
$r.Find.Execute("Experimentos desenvolvidos em laboratório didático, realização de relatórios para cada experimento e de testes sobre o experimento em estudo.", $true, $false, $false, $false, $false, $true, 1, $false, "Média aritmética de duas provas escritas, testes, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + P2 + TR)/3", 2) | Out-Null

# Para 19: "Bibliografia" body
$p = $d.Paragraphs.Item(19)
$p.Range.Text = "6495737 - Durval Rodrigues Junior"
